$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "49.912.98"
$ws.Range("E2").Value = "  +4.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.666.63"
$ws.Range("E3").Value = "  +7.52%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "113.58"
$ws.Range("E5").Value = "  +8.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "326.09"
$ws.Range("E6").Value = "  +2.86%  "
$ws.Range("E7").Value = "  +2.16%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.554"
$ws.Range("E9").Value = "  +3.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.05"
$ws.Range("E10").Value = "  +5.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.06"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0824"
$ws.Range("E12").Value = "  +3.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.35"
$ws.Range("E14").Value = "  +4.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.088.47"
$ws.Range("E15").Value = "  +7.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.664.44"
$ws.Range("E16").Value = "  +6.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.875"
$ws.Range("E17").Value = "  +6.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "49.851.79"
$ws.Range("E18").Value = "  +4.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.13"
$ws.Range("E19").Value = "  +3.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.77"
$ws.Range("E20").Value = "  +4.01%  "
$ws.Range("E21").Value = "  -1.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0959"
$ws.Range("E22").Value = "  +3.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.42"
$ws.Range("E23").Value = "  +2.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "277.50"
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("E25").Value = "  +4.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.84"
$ws.Range("E26").Value = "  +4.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("E28").Value = "  +5.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.143"
$ws.Range("E30").Value = "  +4.22%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.19"
$ws.Range("E31").Value = "  +4.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.27"
$ws.Range("E32").Value = "  +2.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.48"
$ws.Range("E33").Value = "  +4.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.63"
$ws.Range("E34").Value = "  +4.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0817"
$ws.Range("E35").Value = "  +6.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.09"
$ws.Range("E36").Value = "  +12.78%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.07"
$ws.Range("E38").Value = "  +7.42%  "
$ws.Range("E39").Value = "  +10.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "124.87"
$ws.Range("E40").Value = "  +2.24%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.92"
$ws.Range("E41").Value = "  +4.51%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.113"
$ws.Range("E42").Value = "  +2.37%  "
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0319"
$ws.Range("E44").Value = "  +6.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.114.83"
$ws.Range("E45").Value = "  +6.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.30"
$ws.Range("E46").Value = "  +5.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.27"
$ws.Range("E47").Value = "  +14.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.02"
$ws.Range("E48").Value = "  +7.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.04"
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.38"
$ws.Range("E50").Value = "  +5.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.38"
$ws.Range("E51").Value = "  +6.36%  "
